$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: advance the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Step 2: correct the two unit prices that were wrong
$ws.Range("D22").Value = 141
$ws.Range("D37").Value = 196
